# Update feed logs and data lake files
# Append two new log rows (4 and 5) to Sheet1, mirroring the existing
# run_id / rss_url_id / date / response / item_count columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "2024-06-14 20:13:09"
$ws.Range("D4").Value = 200
$ws.Range("E4").Value = 0

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "2024-06-14 20:13:09"
$ws.Range("D5").Value = 200
$ws.Range("E5").Value = 0
